$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Statistik")

# Insert a new column before the current column G ("icdRd_no" is F, "pt_no" is G);
# this shifts "pt_no" (and its data) from G to H.
$ws.Range("G1").EntireColumn.Insert()

# Copy the header formatting (bold, centered) from the neighboring header cell
# so the new header matches the existing header row styling.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New header and value for the extended RD count column
$ws.Range("G1").Value = "icdRd_no_ext"
$ws.Range("G2").Value = 297
